$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text for the "Participant ID" query cell (B2), replacing the old failing query
# that had been producing incorrect / unsorted results (commit: "Failed test cases
# from Instrument model filter"). Using a literal here-string so backticks, single
# quotes and $ characters inside the Cypher text are preserved verbatim.
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.instrument_model in ['Illumina HiSeq 2500']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id Limit 100
'@

$ws.Range("B2").Value = $newParticipantQuery

# Row 2 needs a taller custom height to fit the longer replacement query text.
$ws.Rows.Item(2).RowHeight = 300

# Move the active selection to B2 (it was C2 before), which also clears the
# previously scrolled-away topLeftCell on the sheet view.
$ws.Range("B2").Select()
